# Apply change to row 21: rename model label and update sampled values (E21:X21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "hispnditFilter"

$ws.Range("E21").Value = 0.7010818321629645
$ws.Range("F21").Value = 0.7594578184086344
$ws.Range("G21").Value = 0.9456563717395856
$ws.Range("H21").Value = 0.9519998607969392
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 0.6773292179153978
$ws.Range("K21").Value = 0.9466433269072861
$ws.Range("L21").Value = 0.8467233258443898
$ws.Range("M21").Value = 0.3099058898395918
$ws.Range("N21").Value = 0.9273430168101348
$ws.Range("O21").Value = 0.9666784585856899
$ws.Range("P21").Value = 0.944557729599292
$ws.Range("Q21").Value = 0.8718373137109566
$ws.Range("R21").Value = 0.940175527682934
$ws.Range("S21").Value = 0.943025451569925
$ws.Range("T21").Value = 0.9052271878411812
$ws.Range("U21").Value = 0.3556506702809143
$ws.Range("V21").Value = 0.9359142068168036
$ws.Range("W21").Value = 0.9131209728303235
$ws.Range("X21").Value = 0.9400022208139552
